$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$full = $tr.Text
$target = 'readme.txt'
$idx0 = $full.IndexOf($target)
if ($idx0 -lt 0) {
    throw "Could not locate 'readme.txt' in shape text"
}

# 1-based character position where "readme.txt" begins
$start1 = $idx0 + 1

# Replace the ".txt" portion (4 chars right after "readme") with ".md` "
# (backtick kept literal via single-quoted string)
$dotStart = $start1 + 6
$dotRange = $tr.Characters($dotStart, 4)
$dotRange.Text = '.md` '

# Now remove the stray leading backtick + space that used to precede "file"
$full2 = $tr.Text
$idx2 = $full2.IndexOf('` file')
if ($idx2 -lt 0) {
    throw "Could not locate '\` file' in shape text"
}
$fileRange = $tr.Characters($idx2 + 1, 6)
$fileRange.Text = 'file'
